$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 10-16 (existing rows): new label (col B) + new numeric values (C:P) ---
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.9923753196759012
$ws.Range("D10").Value = 0.9879300160609247
$ws.Range("E10").Value = 0.995896977005419
$ws.Range("F10").Value = 0.9946912582886986
$ws.Range("G10").Value = 0.9923753196759012
$ws.Range("H10").Value = 0.9879300160609247
$ws.Range("I10").Value = 0.9921929612759364
$ws.Range("J10").Value = 0.9948219900607478
$ws.Range("K10").Value = 0.9946551405311592
$ws.Range("L10").Value = 0.9824036829982525
$ws.Range("M10").Value = 0.9923753196759012
$ws.Range("N10").Value = 0.9919134965331718
$ws.Range("O10").Value = 0.9927233927577359
$ws.Range("P10").Value = 0.9918709182371299

$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 1.007824908111867
$ws.Range("D11").Value = 0.9765899659437869
$ws.Range("E11").Value = 0.9928337984043685
$ws.Range("F11").Value = 0.9890144414424805
$ws.Range("G11").Value = 1.007824908111867
$ws.Range("H11").Value = 0.9765899659437869
$ws.Range("I11").Value = 0.9969096996227621
$ws.Range("J11").Value = 0.9885215246904111
$ws.Range("K11").Value = 0.9960292800793002
$ws.Range("L11").Value = 0.9799437275517902
$ws.Range("M11").Value = 1.007824908111867
$ws.Range("N11").Value = 0.9847118821740777
$ws.Range("O11").Value = 0.9915657784756257
$ws.Range("P11").Value = 0.9909584182308458

$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 1.007618925378141
$ws.Range("D12").Value = 0.9768043163281916
$ws.Range("E12").Value = 0.9928101652743769
$ws.Range("F12").Value = 0.9890278428509218
$ws.Range("G12").Value = 1.007618925378141
$ws.Range("H12").Value = 0.9768043163281916
$ws.Range("I12").Value = 0.9968931884101567
$ws.Range("J12").Value = 0.9885020357262566
$ws.Range("K12").Value = 0.9959866672158417
$ws.Range("L12").Value = 0.9801553087160731
$ws.Range("M12").Value = 1.007618925378141
$ws.Range("N12").Value = 0.9848072408012842
$ws.Range("O12").Value = 0.9915653124579078
$ws.Range("P12").Value = 0.990974806237495

$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 1.007820952547674
$ws.Range("D13").Value = 0.9766163577266859
$ws.Range("E13").Value = 0.9928248546740771
$ws.Range("F13").Value = 0.9890314307499328
$ws.Range("G13").Value = 1.007820952547674
$ws.Range("H13").Value = 0.9766163577266859
$ws.Range("I13").Value = 0.9968974565154913
$ws.Range("J13").Value = 0.9885020057967278
$ws.Range("K13").Value = 0.99598943056112
$ws.Range("L13").Value = 0.9799857281903518
$ws.Range("M13").Value = 1.007820952547674
$ws.Range("N13").Value = 0.9847206062003815
$ws.Range("O13").Value = 0.9915733989245925
$ws.Range("P13").Value = 0.9909585270952577

$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.9923000000000011
$ws.Range("D14").Value = 0.9756039999999983
$ws.Range("E14").Value = 1.000008
$ws.Range("F14").Value = 0.9872840000000004
$ws.Range("G14").Value = 0.9923000000000011
$ws.Range("H14").Value = 0.9756039999999983
$ws.Range("I14").Value = 0.9993879999999994
$ws.Range("J14").Value = 0.9918760000000005
$ws.Range("K14").Value = 0.9991759999999991
$ws.Range("L14").Value = 0.9793559999999998
$ws.Range("M14").Value = 0.9923000000000011
$ws.Range("N14").Value = 0.9878059999999993
$ws.Range("O14").Value = 0.9887990000000001
$ws.Range("P14").Value = 0.9906239999999999

$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0.9438874999999985
$ws.Range("E15").Value = 1.01
$ws.Range("F15").Value = 0.98
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0.9438874999999985
$ws.Range("I15").Value = 1.01
$ws.Range("J15").Value = 0.99
$ws.Range("K15").Value = 1.01
$ws.Range("L15").Value = 0.96
$ws.Range("M15").Value = 1
$ws.Range("N15").Value = 0.9769437499999993
$ws.Range("O15").Value = 0.9834718749999997
$ws.Range("P15").Value = 0.9879859374999997

$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.9974985830400008
$ws.Range("D16").Value = 0.964044755353595
$ws.Range("E16").Value = 1.002897146470402
$ws.Range("F16").Value = 0.984886727475201
$ws.Range("G16").Value = 0.9974985830400008
$ws.Range("H16").Value = 0.964044755353595
$ws.Range("I16").Value = 1.002694749593604
$ws.Range("J16").Value = 0.9909551150079978
$ws.Range("K16").Value = 1.002321419673601
$ws.Range("L16").Value = 0.9737400070143994
$ws.Range("M16").Value = 0.9974985830400008
$ws.Range("N16").Value = 0.9834709509119985
$ws.Range("O16").Value = 0.9873318030847996
$ws.Range("P16").Value = 0.9898798129536002

# --- Add new rows 17-19 ---
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9926994209314527
$ws.Range("D17").Value = 0.9926792305866939
$ws.Range("E17").Value = 0.9924176336895202
$ws.Range("F17").Value = 0.9917824504897105
$ws.Range("G17").Value = 0.9926994209314527
$ws.Range("H17").Value = 0.9926792305866939
$ws.Range("I17").Value = 0.9925580205193487
$ws.Range("J17").Value = 0.9925931930861132
$ws.Range("K17").Value = 0.9924170774831098
$ws.Range("L17").Value = 0.9923148522500058
$ws.Range("M17").Value = 0.9926994209314527
$ws.Range("N17").Value = 0.9925484321381071
$ws.Range("O17").Value = 0.9923946839243444
$ws.Range("P17").Value = 0.9924327348794943

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9900025626663325
$ws.Range("D18").Value = 0.9941342646830896
$ws.Range("E18").Value = 0.9926244190266585
$ws.Range("F18").Value = 0.9913563788377007
$ws.Range("G18").Value = 0.9900025626663325
$ws.Range("H18").Value = 0.9941342646830896
$ws.Range("I18").Value = 0.9926404096185981
$ws.Range("J18").Value = 0.9926101554921665
$ws.Range("K18").Value = 0.9922574907241658
$ws.Range("L18").Value = 0.9930435283609363
$ws.Range("M18").Value = 0.9900025626663325
$ws.Range("N18").Value = 0.9933793418548741
$ws.Range("O18").Value = 0.9920294063034454
$ws.Range("P18").Value = 0.992333651176206

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9887333821653637
$ws.Range("D19").Value = 0.997607525701488
$ws.Range("E19").Value = 0.9914569489173536
$ws.Range("F19").Value = 0.9928943056861634
$ws.Range("G19").Value = 0.9887333821653637
$ws.Range("H19").Value = 0.997607525701488
$ws.Range("I19").Value = 0.9905102786960728
$ws.Range("J19").Value = 0.9939504500144281
$ws.Range("K19").Value = 0.9911396292117644
$ws.Range("L19").Value = 0.9963040968158872
$ws.Range("M19").Value = 0.9887333821653637
$ws.Range("N19").Value = 0.9945322373094208
$ws.Range("O19").Value = 0.9926730406175921
$ws.Range("P19").Value = 0.9928245771510651

# Apply the same style as the existing A-column index cells (bold, centered, bordered) to the new rows
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
